# Daily Work Log: freeze yesterday's auto-date row to static values and
# log a new "today" entry (liquid glass buttons + GIT) in row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Prime row 3's formatting by copying it from row 2 (same columns)
#        *before* touching row 2's content, so the new row picks up the
#        existing date / vertical-center styles instead of inventing new
#        style/number-format entries.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C3").PasteSpecial(-4122)   # xlPasteFormats (C3 uses B/D2's plain style, not C2's wrap style)

$excel.CutCopyMode = $false

# --- 2) Fill in today's new row with the live auto-date formulas and the
#        new log entry.
$ws.Range("A3").Formula = "=TODAY()"
$ws.Range("B3").Formula = '=TEXT(TODAY(), "dddd")'
$ws.Range("C3").Value = "Added new liquid glass button features for smooth transision and also connected the project to GIT"
$ws.Range("D3").Value = "2hr"

# --- 3) Freeze row 2 (yesterday) from live formulas into the static
#        values they last displayed.
$ws.Range("A2").Value = 45982
$ws.Range("B2").Value = "Friday"

# --- 4) Match the saved selection.
$ws.Range("C3").Select()
